$wb = $excel.ActiveWorkbook

# Sheet ALC, row 18 (diff hunk 0)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 960.2
$ws.Range("I18").Value = 967
$ws.Range("J18").Value = 950
$ws.Range("K18").Value = 967
$ws.Range("L18").Value = 950
$ws.Range("M18").Value = -683
$ws.Range("N18").Value = -1518

# Sheet ALC, row 107 (diff hunk 1)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 6992.2666
$ws.Range("I107").Value = 8699.083000000001
$ws.Range("J107").Value = 165
$ws.Range("K107").Value = 8699.083000000001
$ws.Range("L107").Value = 165
$ws.Range("M107").Value = -6779.083000000001
$ws.Range("N107").Value = -4005

# Sheet ALC, row 132 (diff hunk 2)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 10828.9795
$ws.Range("I132").Value = 1671.2073
$ws.Range("J132").Value = 55001.766
$ws.Range("K132").Value = 5013.6219
$ws.Range("L132").Value = 165005.298
$ws.Range("M132").Value = -2483.6219
$ws.Range("N132").Value = -170065.298

# Sheet ALC, row 137 (diff hunk 3)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3670.0981
$ws.Range("I137").Value = 1204.375
$ws.Range("J137").Value = 5861.852
$ws.Range("K137").Value = 3613.125
$ws.Range("L137").Value = 17585.556
$ws.Range("M137").Value = -1063.125
$ws.Range("N137").Value = -22685.556

# Sheet ALC, row 138 (diff hunk 4)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1839.41
$ws.Range("I138").Value = 822.4035
$ws.Range("J138").Value = 3187.535
$ws.Range("K138").Value = 2467.2105
$ws.Range("L138").Value = 9562.605
$ws.Range("M138").Value = 2672.7895
$ws.Range("N138").Value = -19842.605

# Sheet ARM, row 4 (diff hunk 5)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 100
$ws.Range("I4").Value = 100
$ws.Range("K4").Value = 100
$ws.Range("M4").Value = 16

# Sheet ARM, row 32 (diff hunk 6)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8990.48
$ws.Range("I32").Value = 8307.182000000001
$ws.Range("J32").Value = 14001.333
$ws.Range("K32").Value = 8307.182000000001
$ws.Range("L32").Value = 14001.333
$ws.Range("M32").Value = -8020.182000000001
$ws.Range("N32").Value = -14575.333

# Sheet ARM, row 37 (diff hunk 7)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 40884.8
$ws.Range("J37").Value = 46796.668
$ws.Range("L37").Value = 46796.668
$ws.Range("N37").Value = -47342.668

# Sheet ARM, row 44 (diff hunk 8)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 34992.332
$ws.Range("J44").Value = 34992.332
$ws.Range("L44").Value = 34992.332
$ws.Range("N44").Value = -35968.332

# Sheet ARM, row 55 (diff hunk 9)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 41952.668
$ws.Range("J55").Value = 41952.668
$ws.Range("L55").Value = 41952.668
$ws.Range("N55").Value = -42582.668

# Sheet ARM, row 82 (diff hunk 10)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H82").Value = 26000
$ws.Range("J82").Value = 26000
$ws.Range("L82").Value = 26000
$ws.Range("N82").Value = -26722

# Sheet ARM, row 85 (diff hunk 11)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H85").Value = 26000
$ws.Range("J85").Value = 26000
$ws.Range("L85").Value = 26000
$ws.Range("N85").Value = -28496

# Sheet ARM, row 102 (diff hunk 12)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 9886.814
$ws.Range("I102").Value = 1755.5333
$ws.Range("J102").Value = 20050.916
$ws.Range("K102").Value = 1755.5333
$ws.Range("L102").Value = 20050.916
$ws.Range("M102").Value = -133.5333000000001
$ws.Range("N102").Value = -23294.916

# Sheet ARM, row 105 (diff hunk 13)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H105").Value = 49362
$ws.Range("J105").Value = 49362
$ws.Range("L105").Value = 49362
$ws.Range("N105").Value = -56350

# Sheet ARM, row 122 (diff hunk 14)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3104.5
$ws.Range("I122").Value = 3262.2856
$ws.Range("K122").Value = 9786.856800000001
$ws.Range("M122").Value = -7336.856800000001

# Sheet BSM, row 99 (diff hunk 15)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2032.25
$ws.Range("I99").Value = 1880.4138
$ws.Range("K99").Value = 1880.4138
$ws.Range("M99").Value = -382.4138

# Sheet BSM, row 107 (diff hunk 16)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1425.38
$ws.Range("I107").Value = 1250.5
$ws.Range("K107").Value = 1250.5
$ws.Range("M107").Value = 669.5

# Sheet CRP, row 4 (diff hunk 17)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 45640704
$ws.Range("I4").Value = 500000260
$ws.Range("J4").Value = 204749.1
$ws.Range("K4").Value = 500000260
$ws.Range("L4").Value = 204749.1
$ws.Range("M4").Value = -500000148
$ws.Range("N4").Value = -204973.1

# Sheet CRP, row 22 (diff hunk 18)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1262.7307
$ws.Range("I22").Value = 510.57895
$ws.Range("J22").Value = 3304.2856
$ws.Range("K22").Value = 510.57895
$ws.Range("L22").Value = 3304.2856
$ws.Range("M22").Value = -160.57895
$ws.Range("N22").Value = -4004.2856

# Sheet CRP, row 31 (diff hunk 19)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2934.7778
$ws.Range("I31").Value = 1164.0714
$ws.Range("J31").Value = 3633.0845
$ws.Range("K31").Value = 1164.0714
$ws.Range("L31").Value = 3633.0845
$ws.Range("M31").Value = -869.0714
$ws.Range("N31").Value = -4223.0845

# Sheet CRP, row 34 (diff hunk 20)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2934.7778
$ws.Range("I34").Value = 1164.0714
$ws.Range("J34").Value = 3633.0845
$ws.Range("K34").Value = 1164.0714
$ws.Range("L34").Value = 3633.0845
$ws.Range("M34").Value = -962.0714
$ws.Range("N34").Value = -4037.0845

# Sheet CRP, row 58 (diff hunk 21)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1266.127
$ws.Range("I58").Value = 1075.7727
$ws.Range("J58").Value = 1706.9474
$ws.Range("K58").Value = 1075.7727
$ws.Range("L58").Value = 1706.9474
$ws.Range("M58").Value = -872.7727
$ws.Range("N58").Value = -2112.9474

# Sheet CRP, row 132 (diff hunk 22)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 46310.594
$ws.Range("I132").Value = 1578
$ws.Range("J132").Value = 144722.3
$ws.Range("K132").Value = 4734
$ws.Range("L132").Value = 434166.9
$ws.Range("M132").Value = -2204
$ws.Range("N132").Value = -439226.9

# Sheet CRP, row 136 (diff hunk 23)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1266.127
$ws.Range("I136").Value = 1075.7727
$ws.Range("J136").Value = 1706.9474
$ws.Range("K136").Value = 3227.3181
$ws.Range("L136").Value = 5120.8422
$ws.Range("M136").Value = -677.3181
$ws.Range("N136").Value = -10220.8422

# Sheet CUL, row 33 (diff hunk 24)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 20534660
$ws.Range("I33").Value = 10
$ws.Range("J33").Value = 22245882
$ws.Range("K33").Value = 60
$ws.Range("L33").Value = 133475292
$ws.Range("M33").Value = 223
$ws.Range("N33").Value = -133475858

# Sheet CUL, row 48 (diff hunk 25)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H48").Value = 500001150
$ws.Range("J48").Value = 500001150
$ws.Range("L48").Value = 1500003450
$ws.Range("N48").Value = -1500003950

# Sheet CUL, row 50 (diff hunk 26)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 83333944
$ws.Range("I50").Value = 400
$ws.Range("J50").Value = 125000710
$ws.Range("K50").Value = 1200
$ws.Range("L50").Value = 375002130
$ws.Range("M50").Value = -719
$ws.Range("N50").Value = -375003092

# Sheet CUL, row 53 (diff hunk 27)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H53").Value = 83333944
$ws.Range("I53").Value = 400
$ws.Range("J53").Value = 125000710
$ws.Range("K53").Value = 1200
$ws.Range("L53").Value = 375002130
$ws.Range("M53").Value = -719
$ws.Range("N53").Value = -375003092

# Sheet CUL, row 88 (diff hunk 28)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 58830076
$ws.Range("J88").Value = 58830076
$ws.Range("L88").Value = 176490228
$ws.Range("N88").Value = -176491084

# Sheet CUL, row 91 (diff hunk 29)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H91").Value = 58830076
$ws.Range("J91").Value = 58830076
$ws.Range("L91").Value = 176490228
$ws.Range("N91").Value = -176493192

# Sheet CUL, row 129 (diff hunk 30)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 116780.195
$ws.Range("I129").Value = 300880.3
$ws.Range("J129").Value = 1717.625
$ws.Range("K129").Value = 902640.8999999999
$ws.Range("L129").Value = 5152.875
$ws.Range("M129").Value = -897640.8999999999
$ws.Range("N129").Value = -15152.875

# Sheet CUL, row 137 (diff hunk 31)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 38471380
$ws.Range("I137").Value = 4426.6665
$ws.Range("J137").Value = 50011464
$ws.Range("K137").Value = 13279.9995
$ws.Range("L137").Value = 150034392
$ws.Range("M137").Value = -8179.999500000002
$ws.Range("N137").Value = -150044592

# Sheet GSM, row 5 (diff hunk 32)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 10000
$ws.Range("I5").Value = 10000
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 10000
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -9888
$ws.Range("N5").ClearContents()

# Sheet GSM, row 39 (diff hunk 33)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 23000
$ws.Range("J39").Value = 23000
$ws.Range("L39").Value = 23000
$ws.Range("N39").Value = -24064

# Sheet GSM, row 126 (diff hunk 34)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 11828.546
$ws.Range("I126").Value = 15214.25
$ws.Range("J126").Value = 2800
$ws.Range("K126").Value = 45642.75
$ws.Range("L126").Value = 8400
$ws.Range("M126").Value = -43172.75
$ws.Range("N126").Value = -13340

# Sheet LTW, row 22 (diff hunk 35)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 521.6
$ws.Range("I22").Value = 429.4
$ws.Range("K22").Value = 429.4
$ws.Range("M22").Value = -134.4

# Sheet LTW, row 27 (diff hunk 36)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 521.6
$ws.Range("I27").Value = 429.4
$ws.Range("K27").Value = 429.4
$ws.Range("M27").Value = -322.4

# Sheet LTW, row 40 (diff hunk 37)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5596.5
$ws.Range("I40").Value = 3148.5
$ws.Range("J40").Value = 10492.5
$ws.Range("K40").Value = 3148.5
$ws.Range("L40").Value = 10492.5
$ws.Range("M40").Value = -3012.5
$ws.Range("N40").Value = -10764.5

# Sheet LTW, row 46 (diff hunk 38)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 10487.625
$ws.Range("J46").Value = 11917
$ws.Range("L46").Value = 11917
$ws.Range("N46").Value = -12293

# Sheet LTW, row 132 (diff hunk 39)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2727.608
$ws.Range("I132").Value = 1840
$ws.Range("J132").Value = 3897.6365
$ws.Range("K132").Value = 5520
$ws.Range("L132").Value = 11692.9095
$ws.Range("M132").Value = -2990
$ws.Range("N132").Value = -16752.9095

# Sheet WVR, row 2 (diff hunk 40)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").ClearContents()

# Sheet WVR, row 51 (diff hunk 41)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 10000
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()

# Sheet WVR, row 62 (diff hunk 42)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3266.65
$ws.Range("I62").Value = 3166.5
$ws.Range("J62").Value = 3277.7778
$ws.Range("K62").Value = 3166.5
$ws.Range("L62").Value = 3277.7778
$ws.Range("M62").Value = -2542.5
$ws.Range("N62").Value = -4525.7778

# Sheet WVR, row 65 (diff hunk 43)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 3266.65
$ws.Range("I65").Value = 3166.5
$ws.Range("J65").Value = 3277.7778
$ws.Range("K65").Value = 15832.5
$ws.Range("L65").Value = 16388.889
$ws.Range("M65").Value = -12712.5
$ws.Range("N65").Value = -22628.889

# Sheet WVR, row 122 (diff hunk 44)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 671.3570999999999
$ws.Range("I122").Value = 649.4545000000001
$ws.Range("K122").Value = 1948.3635
$ws.Range("M122").Value = 501.6364999999998
